$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the hours value for the row dated 2024-05-20 (row 64):
# B64 changes from 5.5 to 6.5. Downstream formulas (C64 running total,
# D2 = SUM(B:B), F2 = 40*D2) will recalculate automatically.
$ws.Range("B64").Value = 6.5

$excel.CalculateFullRebuild()
